# Generate Report for Handback
#
# Refreshes the "Latest HO Xliff Generate Date" / "Correspond Handoff
# Datetime" / "Correspond Handback DateTime" timestamps for the
# 803707bd-1191-4ae0-83a5-dc42dd347716.md file row, across the
# Overview, zh-cn and de-de sheets, reflecting a new handback report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 4 corresponds to 803707bd-1191-4ae0-83a5-dc42dd347716.md
$wsOverview.Range("G4").Value = "2016-11-09 17:40:57"

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 4 corresponds to 803707bd-1191-4ae0-83a5-dc42dd347716.md
$wsZhCn.Range("H4").Value = "2016-11-09 17:40:42"
$wsZhCn.Range("K4").Value = "2016-11-09 17:41:35"

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 4 corresponds to 803707bd-1191-4ae0-83a5-dc42dd347716.md
$wsDeDe.Range("H4").Value = "2016-11-09 17:40:57"
$wsDeDe.Range("K4").Value = "2016-11-09 17:41:54"
